$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1961414790996785
$ws.Range("C2").Value = 0.5434083601286174
$ws.Range("J2").Value = 0.00964630225080386
$ws.Range("P2").Value = 0.135048231511254
$ws.Range("S2").Value = 0.1157556270096463
$ws.Range("B3").Value = 0.005882352941176471
$ws.Range("C3").Value = 0.005882352941176471
$ws.Range("J3").Value = 0.01176470588235294
$ws.Range("P3").Value = 0.7588235294117647
$ws.Range("S3").Value = 0.2176470588235294
$ws.Range("P4").Value = 0.7924528301886793
$ws.Range("S4").Value = 0.2075471698113208
$ws.Range("P5").Value = 0.75
$ws.Range("S5").Value = 0.25
$ws.Range("B6").Value = 0.08695652173913043
$ws.Range("D6").Value = 0.02173913043478261
$ws.Range("E6").Value = 0.01630434782608696
$ws.Range("F6").Value = 0.06521739130434782
$ws.Range("J6").Value = 0.1739130434782609
$ws.Range("O6").Value = 0.005434782608695652
$ws.Range("Q6").Value = 0.1521739130434783
$ws.Range("R6").Value = 0.08695652173913043
$ws.Range("S6").Value = 0.391304347826087
$ws.Range("B7").Value = 0.09444444444444444
$ws.Range("D7").Value = 0.03888888888888889
$ws.Range("F7").Value = 0.08888888888888889
$ws.Range("J7").Value = 0.1277777777777778
$ws.Range("O7").Value = 0.01666666666666667
$ws.Range("Q7").Value = 0.1777777777777778
$ws.Range("R7").Value = 0.05
$ws.Range("S7").Value = 0.4055555555555556
$ws.Range("B8").Value = 0.08955223880597014
$ws.Range("D8").Value = 0.02238805970149254
$ws.Range("E8").Value = 0.002487562189054726
$ws.Range("F8").Value = 0.03980099502487562
$ws.Range("J8").Value = 0.1169154228855721
$ws.Range("O8").Value = 0.01741293532338309
$ws.Range("Q8").Value = 0.1865671641791045
$ws.Range("R8").Value = 0.07960199004975124
$ws.Range("S8").Value = 0.445273631840796
$ws.Range("B9").Value = 0.163265306122449
$ws.Range("D9").Value = 0.04081632653061224
$ws.Range("F9").Value = 0.05102040816326531
$ws.Range("J9").Value = 0.05612244897959184
$ws.Range("O9").Value = 0.01020408163265306
$ws.Range("Q9").Value = 0.1479591836734694
$ws.Range("R9").Value = 0.1326530612244898
$ws.Range("S9").Value = 0.3979591836734694
$ws.Range("B10").Value = 0.1301247771836007
$ws.Range("D10").Value = 0.0249554367201426
$ws.Range("E10").Value = 0.00089126559714795
$ws.Range("F10").Value = 0.07308377896613191
$ws.Range("J10").Value = 0.09090909090909091
$ws.Range("O10").Value = 0.017825311942959
$ws.Range("Q10").Value = 0.2183600713012478
$ws.Range("R10").Value = 0.08021390374331551
$ws.Range("S10").Value = 0.3636363636363636
$ws.Range("G11").Value = 0.1341463414634146
$ws.Range("J11").Value = 0.08536585365853659
$ws.Range("K11").Value = 0.1463414634146341
$ws.Range("L11").Value = 0.6300813008130082
$ws.Range("S11").Value = 0.004065040650406504
$ws.Range("G12").Value = 0.8064516129032258
$ws.Range("J12").Value = 0.167741935483871
$ws.Range("K12").Value = 0.006451612903225806
$ws.Range("L12").Value = 0.006451612903225806
$ws.Range("S12").Value = 0.01290322580645161
$ws.Range("G13").Value = 0.675
$ws.Range("J13").Value = 0.325
$ws.Range("F15").Value = 0.01941747572815534
$ws.Range("H15").Value = 0.2184466019417476
$ws.Range("I15").Value = 0.05825242718446602
$ws.Range("J15").Value = 0.3300970873786408
$ws.Range("K15").Value = 0.06796116504854369
$ws.Range("M15").Value = 0.009708737864077669
$ws.Range("O15").Value = 0.06310679611650485
$ws.Range("S15").Value = 0.2330097087378641
$ws.Range("F16").Value = 0.01449275362318841
$ws.Range("H16").Value = 0.1980676328502415
$ws.Range("I16").Value = 0.07729468599033816
$ws.Range("J16").Value = 0.4057971014492754
$ws.Range("K16").Value = 0.09178743961352658
$ws.Range("M16").Value = 0.00966183574879227
$ws.Range("N16").Value = 0.004830917874396135
$ws.Range("O16").Value = 0.08695652173913043
$ws.Range("S16").Value = 0.1111111111111111
$ws.Range("F17").Value = 0.009852216748768473
$ws.Range("H17").Value = 0.1896551724137931
$ws.Range("I17").Value = 0.1280788177339902
$ws.Range("J17").Value = 0.4014778325123153
$ws.Range("K17").Value = 0.0960591133004926
$ws.Range("M17").Value = 0.01970443349753695
$ws.Range("O17").Value = 0.05911330049261083
$ws.Range("S17").Value = 0.0960591133004926
$ws.Range("F18").Value = 0.005747126436781609
$ws.Range("H18").Value = 0.1494252873563219
$ws.Range("I18").Value = 0.09770114942528736
$ws.Range("J18").Value = 0.4597701149425287
$ws.Range("K18").Value = 0.08620689655172414
$ws.Range("M18").Value = 0.01724137931034483
$ws.Range("O18").Value = 0.05172413793103448
$ws.Range("S18").Value = 0.132183908045977
$ws.Range("F19").Value = 0.01
$ws.Range("H19").Value = 0.2
$ws.Range("I19").Value = 0.09363636363636364
$ws.Range("J19").Value = 0.4145454545454546
$ws.Range("K19").Value = 0.1063636363636364
$ws.Range("M19").Value = 0.02272727272727273
$ws.Range("N19").Value = 0.0009090909090909091
$ws.Range("O19").Value = 0.07818181818181819
$ws.Range("S19").Value = 0.07363636363636364
